# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for zh-cn (B) and de-de (C), and latest handoff date (D)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-31-12 20:31:07"

# zh-cn detail sheet: Status (C) and Latest Handoff Datetime (E)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-12 20:31:03"

# de-de detail sheet: Status (C) and Latest Handoff Datetime (E)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-12 20:31:07"
